$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G ("K") rows 2-22, per the diff
$values = @{
    2  = 12
    3  = 6
    4  = 11
    5  = 6
    6  = 3
    7  = 9
    8  = 4
    9  = 4
    10 = 7
    11 = 6
    12 = 6
    13 = 5
    14 = 6
    15 = 7
    16 = 5
    17 = 8
    18 = 3
    19 = 12
    20 = 3
    21 = 5
    22 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
